# Auto-update script (替自动更新Excel文件 process)
#
# For every data row (row 2 .. last used row) on the active sheet:
#   D = total days ("总天"), E = days remaining ("剩余"), F = start date
#   stored as an integer in yyyyMMdd form ("开始时间").
#
# Each run represents one day passing:
#   - If the remaining-days counter (E) is down to 1, the slot is
#     replenished: E is reset back to the total (D) and the start date
#     (F) is rolled forward by 7 calendar days.
#   - Otherwise E is simply decremented by 1 and F is left untouched.
#   - Rows whose F value is not a well-formed 8 digit yyyyMMdd date are
#     left completely unchanged (defensive - mirrors how the original
#     automation skips malformed rows instead of erroring out).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {

    $D = $ws.Cells.Item($row, 4).Value2
    $E = $ws.Cells.Item($row, 5).Value2
    $F = $ws.Cells.Item($row, 6).Value2

    if ($D -eq $null -or $E -eq $null -or $F -eq $null) {
        continue
    }

    $Fstr = [string]([int]$F)

    if ($Fstr.Length -ne 8) {
        # malformed date value (e.g. "202510929") - skip this row
        continue
    }

    if ([int]$E -eq 1) {
        try {
            $startDate = [datetime]::ParseExact($Fstr, "yyyyMMdd", $null)
        } catch {
            continue
        }
        $newDate = $startDate.AddDays(7)
        $newF = [int]$newDate.ToString("yyyyMMdd")

        $ws.Cells.Item($row, 5).Value2 = [int]$D
        $ws.Cells.Item($row, 6).Value2 = $newF
    } else {
        $ws.Cells.Item($row, 5).Value2 = ([int]$E - 1)
    }
}
